$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "486.50") must be
# forced to remain text so Excel does not reinterpret/renormalize them as
# numeric values (losing formatting like trailing zeros).
$textForceCells = @("D5", "D6", "D7", "D8", "D10", "D11", "D12", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price/volume/coin/link values scraped for this run.
$ws.Range("D2").Value = "56.807.44"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.407.23"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "486.50"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "153.75"
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  +19.44%  "
$ws.Range("D9").Value = "2.430.78"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").Value = "6.24"
$ws.Range("E10").Value = "  +10.06%  "
$ws.Range("D11").Value = "0.0997"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "0.330"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "2.825.93"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").Value = "57.057.21"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "20.57"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "0.0000132"
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("D18").Value = "2.432.18"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "4.64"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").Value = "321.11"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "9.97"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "0.995"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").Value = "5.88"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("D24").Value = "57.89"
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "0.400"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("D27").Value = "0.159"
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("D28").Value = "2.537.26"
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").Value = "  -3.36%  "
$ws.Range("D30").Value = "0.0₃0788"
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "151.03"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "18.67"
$ws.Range("E33").Value = "  +3.39%  "
$ws.Range("D34").Value = "1.52"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").Value = "5.30"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("D36").Value = "3.74"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").Value = "1.13"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").Value = "0.818"
$ws.Range("E38").Value = "  -7.72%  "
$ws.Range("D39").Value = "34.15"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").Value = "1.37"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").Value = "3.51"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  +4.79%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").Value = "274.79"
$ws.Range("E44").Value = "  +4.11%  "
$ws.Range("D45").Value = "0.597"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").Value = "0.0532"
$ws.Range("E46").Value = "  -4.66%  "
$ws.Range("D47").Value = "10.24"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "0.0228"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "4.57"
$ws.Range("E49").Value = "  -5.77%  "
$ws.Range("D50").Value = "17.82"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.882.11"
$ws.Range("E51").Value = "  +2.62%  "

# Restore default (General) styling on the forced-text cells so no stray
# number-format style lingers on them.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
